$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Device 53 -> Device 51 update
$ws.Range("C2").Value = "IND_DAU_51"
$ws.Range("E2").Value = "10.75.58.51"
$ws.Range("F2").Value = 409026540
$ws.Range("G2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# Update selection to F2
$ws.Range("F2").Select()

# Set page setup (orientation/paper size) which generates a pageSetup part
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
